# adding TP & Labo Starter
# Replace placeholder markers ("*" -> "p" and "#" -> "t") in the route02 map
# used-range grid to flag the newly introduced TP & Labo starter cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("route02")

# Rows 3-9: columns C through I change from "*" to "p"
foreach ($r in 3..9) {
    $ws.Range("C$r`:I$r").Value = "p"
}

# Rows 57-61: columns G through M change from "*" to "p"
foreach ($r in 57..61) {
    $ws.Range("G$r`:M$r").Value = "p"
}

# Scattered single/paired cells that change from "#" to "t"
$ws.Range("R12").Value = "t"
$ws.Range("F15:G15").Value = "t"
$ws.Range("S43:T43").Value = "t"
$ws.Range("S47:T47").Value = "t"
$ws.Range("F52:G52").Value = "t"
